$wb = $excel.ActiveWorkbook

# Rename sheets: lowercase language labels -> "Russian - <Capitalized>"
$names = @{
    "adjectives"    = "Russian - Adjectives"
    "adverbs"       = "Russian - Adverbs"
    "conjunctions"  = "Russian - Conjunctions"
    "expressions"   = "Russian - Expressions"
    "nouns"         = "Russian - Nouns"
    "verbs"         = "Russian - Verbs"
}

foreach ($sheet in $wb.Worksheets) {
    $oldName = $sheet.Name
    if ($names.ContainsKey($oldName)) {
        $sheet.Name = $names[$oldName]
    }
}

$wsConjunctions = $wb.Worksheets.Item("Russian - Conjunctions")
$wsVerbs = $wb.Worksheets.Item("Russian - Verbs")

# Verbs sheet: selection moved to G30 (no longer the active tab)
$wsVerbs.Activate()
$wsVerbs.Range("G30").Select()

# Give the Verbs sheet an explicit page setup (portrait, paper size 9 / A4)
$wsVerbs.PageSetup.PaperSize = 9
$wsVerbs.PageSetup.Orientation = 1

# Conjunctions becomes the active / selected tab with selection G12
$wsConjunctions.Activate()
$wsConjunctions.Range("G12").Select()
